# This script applies the update to rows 3-6 of the active worksheet:
#  - Row 3 and Row 4 effectively swap their record content (species/location
#    data), while the Ost/Nord (Q/R) coordinates are updated to new rounded
#    values.
#  - Rows 3-6 have their Ost/Nord (Q/R) coordinates rounded to whole numbers.
#  - Rows 3-6 have the Starttid/Sluttid (Z/AB) cells cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: becomes what used to be row 4's record (rounded coordinates) ---
$ws.Range("A3").Value = 111702281
$ws.Range("B3").Value = 89183
$ws.Range("E3").Value = 3215
$ws.Range("F3").Value = "Rödgul trumpetsvamp"
$ws.Range("G3").Value = "Craterellus lutescens"
$ws.Range("H3").Value = "(Fr.) Fr."
$ws.Range("P3").Value = "Kyrkogården (Kyrkogården), Nrk"
$ws.Range("Q3").Value = 516918
$ws.Range("R3").Value = 6574657

# --- Row 4: becomes what used to be row 3's record (rounded coordinates) ---
$ws.Range("A4").Value = 111701829
$ws.Range("B4").Value = 90687
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = "Fjällig taggsvamp s.str."
$ws.Range("G4").Value = "Sarcodon imbricatus s.str."
$ws.Range("H4").Value = "(L.:Fr.) P.Karst."
$ws.Range("P4").Value = "Myrövägen öster (Myrövägen öster), Nrk"
$ws.Range("Q4").Value = 516895
$ws.Range("R4").Value = 6574639

# --- Row 5: only coordinates rounded ---
$ws.Range("Q5").Value = 516924
$ws.Range("R5").Value = 6574667

# --- Row 6: only coordinates rounded ---
$ws.Range("Q6").Value = 516979
$ws.Range("R6").Value = 6574636

# --- Clear the Starttid (Z) and Sluttid (AB) cells for rows 3-6 ---
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()
